$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update 4th-quarter value for 핸드백 (row 4) from 331 to 0
$ws.Range("G4").Value = 0

# Delete entire row 5 (the 벨트 row), shifting rows up
$ws.Rows(5).Delete()
